$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "참조  https://" (double space) -> "참조: https://" (colon + space)
# typo in the answer texts for rows 2-17, column F.
$rng = $ws.Range("F2:F17")
$rng.Replace("참조  https://", "참조: https://")

# Give the E1 header cell ("exclude_keywords") its own style: same as the
# other header cells but vertically centered instead of top-aligned.
$ws.Range("E1").VerticalAlignment = -4108

# Update the active selection to E1 (previously F1).
$ws.Range("E1").Select()
